$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-8 (columns B = "Wn %" and C = "y (t/cu.m.)") were edited: a new
# reading (1.5 / 1.95) was entered on row 5 and the previous readings that
# used to occupy rows 5-7 were shifted down into rows 6-8 (the old row 8
# reading of 6 / 6.45 drops off since the sheet only has rows 4-8 of data).
$ws.Range("B5").Value2 = 1.5
$ws.Range("C5").Value2 = 1.95

$ws.Range("B6").Value2 = 2
$ws.Range("C6").Value2 = 2.45

$ws.Range("B7").Value2 = 3
$ws.Range("C7").Value2 = 3.45

$ws.Range("B8").Value2 = 4.5
$ws.Range("C8").Value2 = 4.95

# The active selection moved from N9 to C9.
$ws.Range("C9").Select()
